$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A so the existing data (A:D) shifts to (B:E)
$ws.Range("A1").EntireColumn.Insert()

# New header and values for the "name" column
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "OV-2_P2_metagenome"
$ws.Range("A3").Value = "OV-2_P3_metagenome"

# Column widths per diff: new col A ~20.1640625, new col B ~68.5
# (old col C / D widths carried over automatically by the column insert above)
$ws.Columns.Item(1).ColumnWidth = 19.3333333
$ws.Columns.Item(2).ColumnWidth = 67.6666667

# Selection on B3 as recorded in sheetView
$ws.Range("B3").Select()
